$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# C9 (row 9, "LogPathInvalid"):
#   before: "로그 경로가 유효하지 않습니까?"
#   after : "로그 경로가 유효하지 않습니다. 번역값이 달라지면"
# ---------------------------------------------------------------------------
$c9 = $ws.Range("C9")
$c9.Value = "로그 경로가 유효하지 않습니다. 번역값이 달라지면"

$runsC9 = @(
    @{Start=1;  Len=2;  Font="돋움"},
    @{Start=3;  Len=1;  Font="Calibri"},
    @{Start=4;  Len=3;  Font="돋움"},
    @{Start=7;  Len=1;  Font="Calibri"},
    @{Start=8;  Len=4;  Font="돋움"},
    @{Start=12; Len=1;  Font="Calibri"},
    @{Start=13; Len=4;  Font="돋움"},
    @{Start=17; Len=2;  Font="Calibri"},
    @{Start=19; Len=4;  Font="돋움"},
    @{Start=23; Len=1;  Font="Calibri"},
    @{Start=24; Len=4;  Font="돋움"}
)

foreach ($r in $runsC9) {
    $chars = $c9.Characters($r.Start, $r.Len)
    $chars.Font.Name = $r.Font
    $chars.Font.Size = 10
}

# ---------------------------------------------------------------------------
# C8 (row 8, "NoWindowToBeCaptured"):
#   before: "캡쳐 가능 영역이 없습니다."
#   after : "캡쳐 가능 영역이 없습니다. 달라진 xlsx와 같이 코드 올리면?"
# ---------------------------------------------------------------------------
$c8 = $ws.Range("C8")
$c8.Value = "캡쳐 가능 영역이 없습니다. 달라진 xlsx와 같이 코드 올리면?"

$runsC8 = @(
    @{Start=1;  Len=2; Font="돋움"},
    @{Start=3;  Len=1; Font="Calibri"},
    @{Start=4;  Len=2; Font="돋움"},
    @{Start=6;  Len=1; Font="Calibri"},
    @{Start=7;  Len=3; Font="돋움"},
    @{Start=10; Len=1; Font="Calibri"},
    @{Start=11; Len=4; Font="돋움"},
    @{Start=15; Len=2; Font="Calibri"},
    @{Start=17; Len=3; Font="돋움"},
    @{Start=20; Len=5; Font="Calibri"},
    @{Start=25; Len=1; Font="돋움"},
    @{Start=26; Len=1; Font="Calibri"},
    @{Start=27; Len=2; Font="돋움"},
    @{Start=29; Len=1; Font="Calibri"},
    @{Start=30; Len=2; Font="돋움"},
    @{Start=32; Len=1; Font="Calibri"},
    @{Start=33; Len=3; Font="돋움"},
    @{Start=36; Len=1; Font="Calibri"}
)

foreach ($r in $runsC8) {
    $chars = $c8.Characters($r.Start, $r.Len)
    $chars.Font.Name = $r.Font
    $chars.Font.Size = 10
}

# ---------------------------------------------------------------------------
# Selection moves from C10 to C9
# ---------------------------------------------------------------------------
[void]$c9.Select()
